# Update {t,b}values locations; add R17 to BOM
#
# The BOM lists parts grouped by reference designator. A new 0-ohm
# resistor (R17) is inserted right after the R15 resistor row (before the
# DS2411 row), pushing every row from the old "DS2411" row through the
# last row ("X2" / ABRACON-ABS05) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25 (the old "DS2411" row). Excel copies the
# formatting of the row above (row 24, the "5.1M" / R15 row) down into the
# newly inserted row, which is what gives the new G25 cell its s="4" style
# (matching the style already used for the DIGIKEY column on R15's row).
$ws.Rows.Item(25).Insert()

# Populate the new BOM row: Qty, Value, Device, Package, Parts, Description.
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "0k"
$ws.Range("C25").Value = "RESISTOR"
$ws.Range("D25").Value = "0402_RES"
$ws.Range("E25").Value = "R17"
$ws.Range("F25").Value = "Resistor"
# No DIGIKEY part number for this one - leave G25 blank.

# Move the active selection to G25, matching the author's cursor position
# when they saved the workbook.
$ws.Range("G25").Select()
